$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial of 45243 (2023-11-13) for every
# data row (rows 2-70). The update bumps that date by one day to 45244
# (2023-11-14) across all of those rows.
$range = $ws.Range("C2:C70")
$range.Value = 45244
